# collect data thuong phat
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D is "last_edited_time". Rows 5-18 currently carry the timestamp
# "2024-08-03T03:29:00.000Z" (shared with D2:D4, which must stay untouched),
# so we update just D5:D18 to the newer "2024-08-03T03:55:00.000Z".
$ws.Range("D5:D18").Value = "2024-08-03T03:55:00.000Z"

# Rows 19-26 carry the sibling timestamp "2024-08-03T03:28:00.000Z" which
# moves forward to "2024-08-03T03:54:00.000Z".
$ws.Range("D19:D26").Value = "2024-08-03T03:54:00.000Z"
